# Naive forecaster bugfix: drop the oldest forecast row (old row 2, date 39400)
# so every remaining row shifts up by one, then refresh the y_0_forecast (C)
# and y_1_forecast (E) columns with the recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete first data row (date 39400 / 2007); everything below
# shifts up one row, matching the new row 2..52 layout.
$ws.Rows.Item(2).Delete()

# Recomputed "y_0_forecast" values (column C) for rows 2..52.
$cValues = @($null, $null, -0.5037688924316441, -0.5555135891318952, 0.169534172659791, 0.8442071301477228, 1.032338390744236, 1.122475521884692, -0.1937612543835177, -0.578174579726376, -0.8331679621937482, -0.7492845378401558, 0.1273541662098365, 0.2751437421933511, -0.08273351073040391, 0.07468705617190707, -0.1151820594382569, -0.05493014849097255, 0.157394256377752, 0.2820931576894115, 0.2895071529679827, 0.3426151435189873, 0.187607693984293, 0.2343541283920114, 0.3712693419885671, 0.1598952850611068, -0.02097628618118463, -0.009430310228020211, -0.03054415496863694, -0.4923796969465988, -2.657403949513992, -2.657403949513992, -0.2885033948250459, -0.5121403324772844, -0.3096364143617802, -0.3096364143617802, -0.2588455356339781, -0.2454721753057276, -0.1730430455425092, -0.1730430455425092, 0.8644693227634503, 0.7038634017465073, 0.6376744206510576, 0.6376744206510576, 0.001611361207976003, 0.20168190406884, 0.1856341247700399, 0.1856341247700399, 0.08117592553187336, -0.06418790329880686, -0.09450306168263811)

# Recomputed "y_1_forecast" values (column E) for rows 2..52.
$eValues = @(-0.3422723562191532, 0.4944284391569687, -0.467076459743887, -0.4782015746048418, -0.03968684591929561, 1.324233212457782, 0.3765075513336269, 0.7478380109886329, 0.6066710853121382, -0.2445716668737163, -0.3970496740026364, -0.2617076051026235, -0.4865818826308876, -0.100009932057743, -0.01252079199893785, 0.3000376062062493, 0.07348980370169844, 0.1740313431290996, 0.06194937150048041, 0.3390041783450259, 0.1190486358061627, 0.1686730364466316, 0.02358181985058216, 0.2210188332817387, 0.3003887663412641, 0.132465972367557, -0.1339126064348495, -0.09571633453315798, 0.035598638033707, -0.3459257698102514, -2.423328265806446, -1.49562970548649, -0.599348850912329, -1.073589070820447, -0.3018961902350958, -0.1048501255800471, -1.151120647939763, -1.238905350026021, -0.8943276391025989, 0.9692952624595019, 0.7652102000489602, 0.4406734233171727, -0.09128981027868299, 0.1544084105021826, 0.01799217181808199, 0.3613321345859122, 0.4108497965175983, 0.3997355152047577, 0.3069836986764551, -0.08988642825158433, -0.3000102673190841)

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $r = $i + 2
    if ($null -ne $cValues[$i]) {
        $ws.Cells.Item($r, 3).Value = $cValues[$i]
    }
    $ws.Cells.Item($r, 5).Value = $eValues[$i]
}
